$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 37.44676833333333
$ws.Range("H2").Value = 112.340305
$ws.Range("I2").Value = 0.6029245723174422
$ws.Range("J2").Value = 0.6029245723174423
$ws.Range("M2").Value = 3.030775
$ws.Range("N2").Value = 9.092325000000001
$ws.Range("Q2").Value = 113.4927292954583
$ws.Range("R2").Value = 1021.434563659125
$ws.Range("S2").Value = 0.6029245723174422
$ws.Range("T2").Value = 0.6029245723174423

# Row 3
$ws.Range("I3").Value = 0.1838793176915316
$ws.Range("J3").Value = 0.1838793176915316
$ws.Range("M3").Value = 3.030775
$ws.Range("N3").Value = 9.092325000000001
$ws.Range("Q3").Value = 34.612896179675
$ws.Range("R3").Value = 311.516065617075
$ws.Range("S3").Value = 0.1838793176915316
$ws.Range("T3").Value = 0.1838793176915316

# Row 4
$ws.Range("G4").Value = 1.244612333333333
$ws.Range("H4").Value = 3.733837
$ws.Range("I4").Value = 0.02003930892236799
$ws.Range("J4").Value = 0.02003930892236799
$ws.Range("M4").Value = 3.030775
$ws.Range("N4").Value = 9.092325000000001
$ws.Range("Q4").Value = 3.772139944558333
$ws.Range("R4").Value = 33.949259501025
$ws.Range("S4").Value = 0.02003930892236799
$ws.Range("T4").Value = 0.02003930892236799

# Row 5
$ws.Range("G5").Value = 10.03858766666667
$ws.Range("H5").Value = 30.115763
$ws.Range("I5").Value = 0.1616297332180864
$ws.Range("J5").Value = 0.1616297332180864
$ws.Range("M5").Value = 3.030775
$ws.Range("N5").Value = 9.092325000000001
$ws.Range("Q5").Value = 30.42470053544167
$ws.Range("R5").Value = 273.8223048189751
$ws.Range("S5").Value = 0.1616297332180864
$ws.Range("T5").Value = 0.1616297332180864

# Row 6
$ws.Range("G6").Value = 1.327177333333333
$ws.Range("H6").Value = 3.981532
$ws.Range("I6").Value = 0.02136867510078605
$ws.Range("J6").Value = 0.02136867510078605
$ws.Range("M6").Value = 3.030775
$ws.Range("N6").Value = 9.092325000000001
$ws.Range("Q6").Value = 4.022375882433334
$ws.Range("R6").Value = 36.2013829419
$ws.Range("S6").Value = 0.02136867510078605
$ws.Range("T6").Value = 0.02136867510078605

# Row 7
$ws.Range("G7").Value = 0.630923
$ws.Range("H7").Value = 1.892769
$ws.Range("I7").Value = 0.01015839274978569
$ws.Range("J7").Value = 0.01015839274978569
$ws.Range("M7").Value = 3.030775
$ws.Range("N7").Value = 9.092325000000001
$ws.Range("Q7").Value = 1.912185655325
$ws.Range("R7").Value = 17.209670897925
$ws.Range("S7").Value = 0.01015839274978569
$ws.Range("T7").Value = 0.01015839274978569
